$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Kyle Jamieson"

# Force every cell in the target range to be stored as text (matches the
# scraped source data, where even numeric-looking fields like runs/balls/sr
# are plain strings) while keeping the cells at the default (unstyled) format.
$dataRange = $ws.Range("A1:M8")
$dataRange.NumberFormat = "@"

$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$data = @(
  @("10th","Royal Challengers Bangalore","Kyle Jamieson","","11","4","1","1","275.00","Kolkata Knight Riders","Chennai","April 18","RCB won by 38 runs"),
  @("31st","Royal Challengers Bangalore","Kyle Jamieson","run out (Varun)","4","12","0","0","33.33","Kolkata Knight Riders","Abu Dhabi","September 20","KKR won by 9 wickets (with 60 balls remaining)"),
  @("19th","Royal Challengers Bangalore","Kyle Jamieson","run out (Imran Tahir)","16","13","1","1","123.07","Chennai Super Kings","Wankhede","April 25","Super Kings won by 69 runs"),
  @("26th","Royal Challengers Bangalore","Kyle Jamieson","","16","11","1","1","145.45","Punjab Kings","Ahmedabad","April 30","Punjab Kings won by 34 runs"),
  @("6th","Royal Challengers Bangalore","Kyle Jamieson","c Pandey b Holder","12","9","2","0","133.33","Sunrisers Hyderabad","Chennai","April 14","RCB won by 6 runs"),
  @("39th","Royal Challengers Bangalore","Kyle Jamieson","","2","2","0","0","100.00","Mumbai Indians","Dubai (DSC)","September 26","RCB won by 54 runs"),
  @("1st","Royal Challengers Bangalore","Kyle Jamieson","run out (Bumrah)","4","4","0","0","100.00","Mumbai Indians","Chennai","April 09","RCB won by 2 wickets")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        $cell = $ws.Cells.Item($r + 2, $c + 1)
        if ($val -eq "") {
            # A literal Value = "" clears the cell instead of storing an
            # empty string, so use an empty-string formula to get a real
            # (non-blank) text cell, matching the scraped "not out" rows.
            $cell.Formula = "=""""" 
        } else {
            $cell.Value = $val
        }
    }
}

# Drop the temporary "@" number format so cells keep the workbook's default
# style (the source file has no explicit per-cell formatting).
$dataRange.ClearFormats()
